$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.002")
# are preserved exactly as text instead of being parsed into floats.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.948.25'
$ws.Range('E2').Value = '  -4.67%  '
$ws.Range('D3').Value = '1.739.09'
$ws.Range('E3').Value = '  -5.11%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '226.51'
$ws.Range('E5').Value = '  -3.90%  '
$ws.Range('D6').Value = '0.5799'
$ws.Range('E6').Value = '  -4.00%  '
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.2728'
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('D9').Value = '23.28'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').Value = '0.06612'
$ws.Range('E10').Value = '  -5.49%  '
$ws.Range('D11').Value = '0.07559'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').Value = '1.737.26'
$ws.Range('E12').Value = '  -5.37%  '
$ws.Range('D13').Value = '4.701'
$ws.Range('E13').Value = '  -1.37%  '
$ws.Range('D14').Value = '0.6030'
$ws.Range('E14').Value = '  -4.77%  '
$ws.Range('D15').Value = '1.976.33'
$ws.Range('E15').Value = '  -5.10%  '
$ws.Range('D16').Value = '74.60'
$ws.Range('E16').Value = '  -4.45%  '
$ws.Range('D17').Value = '0.000008695'
$ws.Range('E17').Value = '  -12.08%  '
$ws.Range('D18').Value = '27.953.66'
$ws.Range('E18').Value = '  -3.62%  '
$ws.Range('D19').Value = '5.328'
$ws.Range('E19').Value = '  -5.14%  '
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('D21').Value = '205.60'
$ws.Range('E21').Value = '  -5.94%  '
$ws.Range('D22').Value = '11.28'
$ws.Range('E22').Value = '  -2.93%  '
$ws.Range('D23').Value = '6.620'
$ws.Range('E23').Value = '  -4.56%  '
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('D25').Value = '150.07'
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('D26').Value = '8.100'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').Value = '0.1233'
$ws.Range('E27').Value = '  -4.80%  '
$ws.Range('D28').Value = '16.13'
$ws.Range('E28').Value = '  -2.70%  '
$ws.Range('D29').Value = '1.384'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('D30').Value = '0.06151'
$ws.Range('E30').Value = '  -4.69%  '
$ws.Range('D31').Value = '1.391'
$ws.Range('E31').Value = '  -3.77%  '
$ws.Range('D32').Value = '3.740'
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('D33').Value = '3.722'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').Value = '1.666'
$ws.Range('E34').Value = '  -4.04%  '
$ws.Range('D35').Value = '1.035'
$ws.Range('E35').Value = '  -5.78%  '
$ws.Range('D36').Value = '0.6415'
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('D37').Value = '2.419'
$ws.Range('E37').Value = '  -4.88%  '
$ws.Range('E38').Value = '  -1.14%  '
$ws.Range('D39').Value = '0.01668'
$ws.Range('E39').Value = '  -5.13%  '
$ws.Range('D40').Value = '1.132.35'
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').Value = '6.175'
$ws.Range('E41').Value = '  -6.48%  '
$ws.Range('D42').Value = '0.8754'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '99.67'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').Value = '1.889.82'
$ws.Range('E45').Value = '  -5.31%  '
$ws.Range('D46').Value = '59.37'
$ws.Range('E46').Value = '  -4.89%  '
$ws.Range('D47').Value = '1.578'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').Value = '0.00000000107'
$ws.Range('E48').Value = '  -5.15%  '
$ws.Range('D49').Value = '8.251'
$ws.Range('E49').Value = '  -3.43%  '
$ws.Range('D50').Value = '0.05377'
$ws.Range('E50').Value = '  -2.23%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.4417'
$ws.Range('E51').Value = '  -2.88%  '
